$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 886 (shifts existing rows 886..927 down to 887..928)
$ws.Rows.Item(886).Insert()

# Force column A to be stored as literal text so the date-like string
# is not auto-converted into a date serial number.
$ws.Range("A886").NumberFormat = "@"
$ws.Range("A886").Value = "2026/02/27"
$ws.Range("B886").Value = "金"
$ws.Range("C886").Value = 16
$ws.Range("D886").Value = 201
